$d = $word.ActiveDocument

$null = $d.Content.Find.Execute("Format has been corrected not the timing", $true, $false, $false, $false, $false, $true, 1, $false, "Umbizo limesahihishwa sio wakati", 2)
$null = $d.Content.Find.Execute("I added 25 seconds to each timing to correct for the intro song -john argentino", $true, $false, $false, $false, $false, $true, 1, $false, "Niliongeza sekunde 25 kwa kila muda ili kusahihisha wimbo wa utangulizi -john argentino", 2)
$null = $d.Content.Find.Execute("The airport problem - subtitles:", $true, $false, $false, $false, $false, $true, 1, $false, "Tatizo la uwanja wa ndege - manukuu:", 2)
$null = $d.Content.Find.Execute("The administrations of three", $true, $false, $false, $false, $false, $true, 1, $false, "Utawala wa tatu", 2)
$null = $d.Content.Find.Execute("neighboring cities: A, B and C decided", $true, $false, $false, $false, $false, $true, 1, $false, "miji jirani: A, B na C waliamua", 2)
$null = $d.Content.Find.Execute("to build an airport dividing the costs of", $true, $false, $false, $false, $false, $true, 1, $false, "kujenga uwanja wa ndege unaogawanya gharama za", 2)
$null = $d.Content.Find.Execute("implementation. The condition on the", $true, $false, $false, $false, $false, $true, 1, $false, "utekelezaji. Hali juu ya", 2)
$null = $d.Content.Find.Execute("choice of the most suitable place is", $true, $false, $false, $false, $false, $true, 1, $false, "uchaguzi wa mahali pa kufaa zaidi ni", 2)
$null = $d.Content.Find.Execute("that the sum of the distances from each", $true, $false, $false, $false, $false, $true, 1, $false, "kwamba jumla ya umbali kutoka kwa kila mmoja", 2)
$null = $d.Content.Find.Execute("city to the airport is as small as", $true, $false, $false, $false, $false, $true, 1, $false, "mji kwa uwanja wa ndege ni ndogo kama", 2)
$null = $d.Content.Find.Execute("possible. The team of experts in charge", $true, $false, $false, $false, $false, $true, 1, $false, "inawezekana. Timu ya wataalam wanaohusika", 2)
$null = $d.Content.Find.Execute("of the work has created a model to get", $true, $false, $false, $false, $false, $true, 1, $false, "ya kazi imeunda mfano wa kupata", 2)
$null = $d.Content.Find.Execute("a preliminary idea of where to place the", $true, $false, $false, $false, $false, $true, 1, $false, "wazo la awali la mahali pa kuweka", 2)
$null = $d.Content.Find.Execute("structure. At their disposal there are", $true, $false, $false, $false, $false, $true, 1, $false, "muundo. Ovyo wao wapo", 2)
$null = $d.Content.Find.Execute("some snails a big metal ring and a long", $true, $false, $false, $false, $false, $true, 1, $false, "konokono wengine pete kubwa ya chuma na ndefu", 2)
$null = $d.Content.Find.Execute("string.", $true, $false, $false, $false, $false, $true, 1, $false, "kamba.", 2)
$null = $d.Content.Find.Execute("Explain how the team can manage to use", $true, $false, $false, $false, $false, $true, 1, $false, "Eleza jinsi timu inaweza kusimamia matumizi", 2)
$null = $d.Content.Find.Execute("the materials to tell approximately the", $true, $false, $false, $false, $false, $true, 1, $false, "nyenzo za kusema takriban", 2)
$null = $d.Content.Find.Execute("ideal location of the airport. Imagine", $true, $false, $false, $false, $false, $true, 1, $false, "eneo bora la uwanja wa ndege. Fikiria", 2)
$null = $d.Content.Find.Execute("that the cities are placed at the", $true, $false, $false, $false, $false, $true, 1, $false, "kwamba miji imewekwa kwenye", 2)
$null = $d.Content.Find.Execute("vertices of a triangle which is", $true, $false, $false, $false, $false, $true, 1, $false, "vipeo vya pembetatu ambayo ni", 2)
$null = $d.Content.Find.Execute("obviously reproduced in scale as", $true, $false, $false, $false, $false, $true, 1, $false, "kwa hakika imetolewa tena kwa kiwango kama", 2)
$null = $d.Content.Find.Execute("shown in figure. This is one possible", $true, $false, $false, $false, $false, $true, 1, $false, "inavyoonyeshwa kwenye takwimu. Hili ni moja linalowezekana", 2)
$null = $d.Content.Find.Execute("setting the rope starts from one nail,", $true, $false, $false, $false, $false, $true, 1, $false, "kuweka kamba huanza kutoka msumari mmoja,", 2)
$null = $d.Content.Find.Execute("goes inside the ring, goes around the", $true, $false, $false, $false, $false, $true, 1, $false, "huenda ndani ya pete, huzunguka", 2)
$null = $d.Content.Find.Execute("other nail, the third nail, inside the", $true, $false, $false, $false, $false, $true, 1, $false, "msumari mwingine, msumari wa tatu, ndani ya", 2)
$null = $d.Content.Find.Execute("ring again and now you can just pull the", $true, $false, $false, $false, $false, $true, 1, $false, "pete tena na sasa unaweza kuvuta tu", 2)
$null = $d.Content.Find.Execute("rope in order to find the point that", $true, $false, $false, $false, $false, $true, 1, $false, "kamba ili kupata uhakika huo", 2)
$null = $d.Content.Find.Execute("you're looking for. In order to reach the", $true, $false, $false, $false, $false, $true, 1, $false, "unatafuta. Ili kufikia", 2)
$null = $d.Content.Find.Execute("point, we have to move the rope a bit", $true, $false, $false, $false, $false, $true, 1, $false, "uhakika, tunapaswa kusonga kamba kidogo", 2)
$null = $d.Content.Find.Execute("because there is some ", $true, $false, $false, $false, $false, $true, 1, $false, "kwa sababu kuna ", 2)
$null = $d.Content.Find.Execute("resistance", $true, $false, $false, $false, $false, $true, 1, $false, "upinzani", 2)
$null = $d.Content.Find.Execute(" caused", $true, $false, $false, $false, $false, $true, 1, $false, " uliosababishwa", 2)
$null = $d.Content.Find.Execute("by the materials that we are using but", $true, $false, $false, $false, $false, $true, 1, $false, "kwa nyenzo ambazo tunatumia lakini", 2)
$null = $d.Content.Find.Execute("after a while you'll reach a position from", $true, $false, $false, $false, $false, $true, 1, $false, "baada ya muda utafikia nafasi kutoka", 2)
$null = $d.Content.Find.Execute("which the ring doesn't move anymore,", $true, $false, $false, $false, $false, $true, 1, $false, "ambayo pete haisogei tena,", 2)
$null = $d.Content.Find.Execute("which is more or less this one. And as", $true, $false, $false, $false, $false, $true, 1, $false, "ambayo ni zaidi au chini ya hii. Na kama", 2)
$null = $d.Content.Find.Execute("between the ring and the nails are", $true, $false, $false, $false, $false, $true, 1, $false, "kati ya pete na misumari ni", 2)
$null = $d.Content.Find.Execute("placed more or less 120 degrees from one", $true, $false, $false, $false, $false, $true, 1, $false, "kuwekwa zaidi au chini ya digrii 120 kutoka kwa moja", 2)
$null = $d.Content.Find.Execute("another which is 1/3 of a circumference,", $true, $false, $false, $false, $false, $true, 1, $false, "nyingine ambayo ni 1/3 ya mduara,", 2)
$null = $d.Content.Find.Execute("and that's the point that we're looking", $true, $false, $false, $false, $false, $true, 1, $false, "na hiyo ndiyo hatua tunayoiangalia", 2)
$null = $d.Content.Find.Execute("for: the minimum distance between the", $true, $false, $false, $false, $false, $true, 1, $false, "kwa: umbali wa chini kati ya", 2)
$null = $d.Content.Find.Execute("nails and the airport when you sum it", $true, $false, $false, $false, $false, $true, 1, $false, "misumari na uwanja wa ndege unapojumlisha", 2)
$null = $d.Content.Find.Execute("ogether", $true, $false, $false, $false, $false, $true, 1, $false, "pamoja", 2)
$null = $d.Content.Find.Execute("[Music]", $true, $false, $false, $false, $false, $true, 1, $false, "[Muziki]", 2)
